# Generate Report for Handback
#
# The localization status report is refreshed once handback (localized
# files flowing back from the translators and being in sync with en-US
# again) has happened for the two source docs tracked in this workbook
# (24222ef3-... and 8c9297c7-...), for both locales (zh-cn, de-de).
#
# This updates:
#  - the "Status" text (Overview + per-locale sheets) from the handoff-time
#    placeholder to a handed-back confirmation,
#  - the newly-populated "Latest Target File" / "Latest Handback File"
#    columns on the per-locale sheets (target doc hyperlink + handback
#    xliff file name), and
#  - "Latest Handback DateTime" with the actual handback timestamps.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$doc1 = "24222ef3-212d-4d3e-83a6-75c7416c571f"
$doc2 = "8c9297c7-c14f-4444-b02e-345cf21e377c"

$doc1Md = "$doc1.md"
$doc2Md = "$doc2.md"

$doc1Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/002bbc19a7f431b1b380d444506eb7928b142f29/e2e/$doc1Md"
$doc2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/002bbc19a7f431b1b380d444506eb7928b142f29/e2e/$doc2Md"

# The status text is noticeably longer than the handoff-time placeholder, so
# the status/"Latest Target File"/"Latest Handback File" columns need to
# widen to comfortably fit it (mirrors the wider columns the report
# generator writes once it has real handback data to show).
$wideColumnWidth = 30 - (5 / 6)
$maxColumnWidth = 40 - (5 / 6)

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-locale status cells.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value2 = $statusText
$overview.Range("F2").Value2 = $statusText
$overview.Range("E3").Value2 = $statusText
$overview.Range("F3").Value2 = $statusText
$overview.Columns.Item(5).ColumnWidth = $wideColumnWidth
$overview.Columns.Item(6).ColumnWidth = $wideColumnWidth

# ---------------------------------------------------------------------
# zh-cn sheet: status text + newly-available target/handback info.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value2 = $statusText
$zhcn.Range("C3").Value2 = $statusText

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $doc1Url, "", "", $doc1Md)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $doc2Url, "", "", $doc2Md)

$zhcn.Range("J2").Value2 = "$doc1.a7d4875fc4555919c4d354939c6f6863e3b47f77.zh-cn.xlf"
$zhcn.Range("J3").Value2 = "$doc2.b3e3e337c85a8ed4ffac93a40cecc2367efb44fc.zh-cn.xlf"

$zhcn.Range("K2").Value2 = "2016-08-18 07:00:59"
$zhcn.Range("K3").Value2 = "2016-08-18 07:00:59"

$zhcn.Columns.Item(3).ColumnWidth = $wideColumnWidth
$zhcn.Columns.Item(9).ColumnWidth = $maxColumnWidth
$zhcn.Columns.Item(10).ColumnWidth = $maxColumnWidth

# ---------------------------------------------------------------------
# de-de sheet: status text + newly-available target/handback info.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value2 = $statusText
$dede.Range("C3").Value2 = $statusText

$dede.Hyperlinks.Add($dede.Range("I2"), $doc1Url, "", "", $doc1Md)
$dede.Hyperlinks.Add($dede.Range("I3"), $doc2Url, "", "", $doc2Md)

$dede.Range("J2").Value2 = "$doc1.a7d4875fc4555919c4d354939c6f6863e3b47f77.de-de.xlf"
$dede.Range("J3").Value2 = "$doc2.b3e3e337c85a8ed4ffac93a40cecc2367efb44fc.de-de.xlf"

$dede.Range("K2").Value2 = "2016-08-18 07:01:19"
$dede.Range("K3").Value2 = "2016-08-18 07:01:19"

$dede.Columns.Item(3).ColumnWidth = $wideColumnWidth
$dede.Columns.Item(9).ColumnWidth = $maxColumnWidth
$dede.Columns.Item(10).ColumnWidth = $maxColumnWidth
